$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 20: was a placeholder duplicate of "Front-end: Edit profile"; now
# becomes the real "Front-end: Dashboard view data management" task, with
# its Finish/Actual Finish dates filled in and an updated comment.
# ---------------------------------------------------------------------------
$ws.Rows.Item(20).RowHeight = 75

$ws.Range("B20").Value = "Front-end: Dashboard view data management"
$ws.Range("D20").Value = 3

$ws.Range("E20").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("F20").Value = 44287

$ws.Range("E20").Copy()
$ws.Range("G20").PasteSpecial(-4122)
$ws.Range("G20").Value = 44287

$ws.Range("I20").Value = "By moving data out of the individual card components, potentially less API calls will be required"

# ---------------------------------------------------------------------------
# Row 21: used to be an almost-blank spacer row; now holds the
# "Front-end: Integrate Threejs scene" task (moved here from row 28) with
# full details added.
# ---------------------------------------------------------------------------
$ws.Rows.Item(21).RowHeight = 120

$ws.Range("A20").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A21").Value = 20

$ws.Range("B21").Value = "Front-end: Integrate Threejs scene"

$ws.Range("C20").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$ws.Range("C21").Value = "LMS v3"

$ws.Range("D20").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("D21").Value = 8

$ws.Range("E21").Value = 44291
$ws.Range("F21").Value = 44300

$ws.Range("I20").Copy()
$ws.Range("I21").PasteSpecial(-4122)
$ws.Range("I21").Value = "Work on this at least till new pages or content added to Adobe XD link`nExpected to take a long time as it was originally not written in React, and I foresee some issues arising from that"

# ---------------------------------------------------------------------------
# Row 28: the Threejs task (and its stray "Implement existing pages first as
# discussed with Sanjay" note) moved up to row 21, so this row goes back to
# being just an empty placeholder row.
# ---------------------------------------------------------------------------
$ws.Range("B28").ClearContents()
$ws.Range("C28").ClearContents()
$ws.Range("H28").ClearContents()

# ---------------------------------------------------------------------------
# Sheet view: scrolled down, zoomed to 100%, and a new active selection.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$win.Zoom = 100
$ws.Range("G21").Select()
